$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("purpose") rows 2 through 22 contain "fullRNASEQ" which should
# be corrected to "fullRNASeq" (fixing formatting/casing on the purpose column).
for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
